# Updated capital structure database
# Apply updated financial metrics for the Kyrgyzstan precious metals rows (row 2 and row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = -0.182089552238806
    $ws.Range("H$r").Value = -0.182089552238806
    $ws.Range("I$r").Value = -0.2194029850746269
    $ws.Range("J$r").Value = -0.2194029850746269
    $ws.Range("K$r").Value = -30.6
    $ws.Range("L$r").Value = -0.4567164179104478

    $ws.Range("U$r").Value = 5.93
    $ws.Range("V$r").Value = 0.03633578431372549
    $ws.Range("W$r").Value = -1.854545454545455
    $ws.Range("X$r").Value = 0.1310963147421401
    $ws.Range("Y$r").Value = -1.985641769287595
    $ws.Range("Z$r").Value = 0.7045956462298875
    $ws.Range("AA$r").Value = -0.1545903880534231
    $ws.Range("AB$r").Value = 0.1114338863607346
    $ws.Range("AC$r").Value = -0.2660242744141577
    $ws.Range("AD$r").Value = 72.2
    $ws.Range("AE$r").Value = 0
    $ws.Range("AF$r").Value = 72.2
    $ws.Range("AG$r").Value = 66.27000000000001
    $ws.Range("AH$r").Value = 0.3067119796091759
    $ws.Range("AI$r").Value = 0.7141444114737884
    $ws.Range("AJ$r").Value = 0.2887959210354295
    $ws.Range("AK$r").Value = 0.6963328780077755
    $ws.Range("AL$r").Value = 12.5
    $ws.Range("AM$r").Value = 12.486
    $ws.Range("AN$r").Value = -6.278260869565218
    $ws.Range("AO$r").Value = -1.176
    $ws.Range("AP$r").Value = -5.762608695652175
    $ws.Range("AQ$r").Value = -1.177318596828448
}
